# Insert a new row above row 157 (pushes existing rows 157.. down by one)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(157).Insert()

# Populate the newly inserted row 157 with the same record as the (now
# shifted-down) row that follows it, but with its own date (44572).
$ws.Range("A157").Value = 3
$ws.Range("B157").Value = "Femacal de La Calera"
$ws.Range("C157").Value = "Coquimbo"
$ws.Range("D157").Value = 44572
$ws.Range("E157").Value = 5
$ws.Range("F157").Value = 100112039
$ws.Range("G157").Value = "Ciboulette"
$ws.Range("H157").Value = "Sin especificar"
$ws.Range("I157").Value = "Primera"
$ws.Range("J157").Value = 160
$ws.Range("K157").Value = 1500
$ws.Range("L157").Value = 1500
$ws.Range("M157").Value = 1500
$ws.Range("N157").Value = "$/docena de atados"
$ws.Range("O157").Value = "Provincia de Quillota"
$ws.Range("P157").Value = 500
$ws.Range("Q157").Value = 3
$ws.Range("R157").Value = "Hortaliza"
